# Actualización desde MV -datos-
# Adds 5 new daily rows (28-09-2021 .. 04-10-2021) to the end of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append starting at row 191
$newRows = @(
    @{ Row = 191; Fecha = "28-09-2021"; B = 29164; C = 5357; D = 6903; E = -27603 },
    @{ Row = 192; Fecha = "29-09-2021"; B = 29164; C = 5357; D = 6903; E = -27132 },
    @{ Row = 193; Fecha = "30-09-2021"; B = 29164; C = 5357; D = 6331; E = -28020 },
    @{ Row = 194; Fecha = "01-10-2021"; B = 29164; C = 5357; D = 6331; E = -27368 },
    @{ Row = 195; Fecha = "04-10-2021"; B = 29164; C = 5357; D = 6331; E = -27016 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $cellA = $ws.Cells.Item($row, 1)

    # Some "DD-MM-YYYY" labels (where both parts are <=12) are ambiguous and
    # would otherwise be auto-recognized as dates. Force them to be stored as
    # plain text (matching the rest of column A), then restore the default
    # "Normal" cell style so no visible formatting change is introduced.
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Fecha
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
}
